$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.212.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.811.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.57%  "

$ws.Range("E4").Value = "  -0.70%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3930"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3494"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.23"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.178"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07552"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.525"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.809.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.170"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.86%  "

$ws.Range("E17").Value = "  +1.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06707"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("E21").Value = "  +2.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.564"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.200.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.403"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.479"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.526"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.015.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "136.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.235"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.018"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08862"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02437"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6942"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.459"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06529"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.607"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2218"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.267"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.530"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6429"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9996"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.875"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.151"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "131.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07206"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.77%  "
